$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: find the LAST paragraph whose range text equals $text (ignoring the
# trailing paragraph mark), scanning from the end of the document backwards.
# ---------------------------------------------------------------------------
function Find-LastParagraph($doc, [string]$text) {
    $count = $doc.Paragraphs.Count
    for ($i = $count; $i -ge 1; $i--) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text
        if ($t.Length -gt 0) {
            $trimmed = $t.Substring(0, $t.Length - 1)
        } else {
            $trimmed = $t
        }
        if ($trimmed -eq $text) {
            return $p
        }
    }
    return $null
}

$xmlHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$xmlFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------------
# Change 1: the trailing "implemented:" paragraph of the "Log 2013-05-06"
# entry is missing an explicit paragraph-mark language tag. Rewrite that
# paragraph's XML so the <w:pPr> carries <w:rPr><w:lang w:val="en-US"/></w:rPr>,
# matching the same paragraph earlier in the document.
# ---------------------------------------------------------------------------
$implementedParagraph = Find-LastParagraph $d "implemented:"
if ($implementedParagraph -ne $null) {
    $implementedXml = $xmlHeader + '<w:p><w:pPr><w:pStyle w:val="Listeafsnit"/><w:ind w:left="0"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>implemented</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>:</w:t></w:r></w:p>' + $xmlFooter
    $implementedParagraph.Range.InsertXML($implementedXml)
}

# ---------------------------------------------------------------------------
# Change 2: append a brand-new "Log 2013-05-09" entry right after the
# "Remove black jack bug" bullet (end of the "Log 2013-05-06" entry) and
# before the document's final (empty) paragraph.
# ---------------------------------------------------------------------------
$lastBugParagraph = Find-LastParagraph $d "Remove black jack bug"
if ($lastBugParagraph -ne $null) {
    $insertRange = $lastBugParagraph.Range
    $insertRange.Collapse(0)

    $newEntryBody = '<w:p><w:pPr><w:pStyle w:val="DefaultStyle"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:b/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>Log 2013-05-09</w:t></w:r></w:p>' +
        '<w:p><w:pPr><w:pStyle w:val="DefaultStyle"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>Driver: Simon</w:t></w:r></w:p>' +
        '<w:p><w:pPr><w:pStyle w:val="DefaultStyle"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Navigator: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>Thelle</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>' +
        '<w:p><w:pPr><w:pStyle w:val="DefaultStyle"/><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">What did we </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>do:</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p>' +
        '<w:p><w:pPr><w:pStyle w:val="Listeafsnit"/><w:ind w:left="0"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>implemented</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>:</w:t></w:r></w:p>' +
        '<w:p><w:pPr><w:pStyle w:val="Listeafsnit"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>Squashed Black jack bug</w:t></w:r></w:p>' +
        '<w:p><w:pPr><w:pStyle w:val="Listeafsnit"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Changed description of death by gun</w:t></w:r></w:p>' +
        '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Plans for next time:</w:t></w:r></w:p>' +
        '<w:p><w:pPr><w:pStyle w:val="Listeafsnit"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Nothing </w:t></w:r></w:p>'

    $newEntryXml = $xmlHeader + $newEntryBody + $xmlFooter
    $insertRange.InsertXML($newEntryXml)
}
